$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update KODE_REKSADANA (column N, row 2) from RD00014 to RD00015
$ws.Range("N2").Value = "RD00015"

# Update PREPARATION (column F, row 2) text to reflect new Kode Reksadana
$ws.Range("F2").Value = "Username : 31246;`nPassword : bni1234;`nRole : 18/19/20/21 - Pimpinan Kelompok Investasi/Pengelolan Investasi/Analis;`nKode Reksadana : RD00015"

# Update the selection to G2
$ws.Range("G2").Select()
